$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.703.40"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "1.704.98"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Formula = "'309.05"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Formula = "'0.9985"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Formula = "'0.3745"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Formula = "'49.29"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("D9").Formula = "'0.3442"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Formula = "'1.191"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Formula = "'0.07474"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Formula = "'1.000"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Formula = "'20.94"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").Formula = "'6.248"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Formula = "'6.959"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "1.706.95"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Formula = "'0.00001129"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Formula = "'0.9989"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Formula = "'84.37"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").Formula = "'17.15"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Formula = "'6.339"
$ws.Range("D23").Formula = "'13.06"
$ws.Range("E23").Value = "  +8.70%  "
$ws.Range("D24").Value = "24.700.55"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").Formula = "'2.426"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Formula = "'2.765"
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").Formula = "'20.21"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Formula = "'150.78"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Formula = "'131.36"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").Value = "1.894.09"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Formula = "'1.187"
$ws.Range("E31").Value = "  +21.64%  "
$ws.Range("D32").Formula = "'6.800"
$ws.Range("E32").Value = "  +6.62%  "
$ws.Range("D33").Formula = "'4.173"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").Formula = "'1.802"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").Formula = "'0.08840"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("D36").Formula = "'13.69"
$ws.Range("E36").Value = "  +10.93%  "
$ws.Range("D37").Formula = "'5.539"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("D38").Formula = "'0.06590"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").Formula = "'8.989"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Formula = "'0.02393"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").Formula = "'0.2232"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("D42").Formula = "'1.276"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Formula = "'0.6447"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Formula = "'13.95"
$ws.Range("E45").Value = "  +6.16%  "
$ws.Range("D46").Formula = "'0.6125"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").Formula = "'3.813"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Formula = "'2.121"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("D49").Formula = "'129.79"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").Formula = "'79.38"
$ws.Range("E51").Value = "  +4.36%  "

Write-Host "Updated crypto prices"